# Auto-generated edit script: updates column F (想去人数 / want-to-go count) values
# across sheets 展览 (sheet1), 演出 (sheet2), and 全部类型 (sheet4),
# matching the target diff. No changes needed on 本地生活 (sheet3).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 69
$ws1.Range("F4").Value = 27
$ws1.Range("F5").Value = 1159
$ws1.Range("F6").Value = 9014
$ws1.Range("F7").Value = 147
$ws1.Range("F8").Value = 244
$ws1.Range("F9").Value = 7144
$ws1.Range("F10").Value = 185
$ws1.Range("F11").Value = 322
$ws1.Range("F12").Value = 5494
$ws1.Range("F13").Value = 10
$ws1.Range("F14").Value = 72
$ws1.Range("F15").Value = 6267
$ws1.Range("F16").Value = 1096
$ws1.Range("F19").Value = 275
$ws1.Range("F20").Value = 148
$ws1.Range("F22").Value = 161
$ws1.Range("F23").Value = 104
$ws1.Range("F24").Value = 10152
$ws1.Range("F25").Value = 84
$ws1.Range("F26").Value = 1929
$ws1.Range("F27").Value = 1973
$ws1.Range("F28").Value = 47
$ws1.Range("F29").Value = 37
$ws1.Range("F30").Value = 2130
$ws1.Range("F31").Value = 0
$ws1.Range("F33").Value = 169
$ws1.Range("F34").Value = 1033
$ws1.Range("F36").Value = 2087
$ws1.Range("F37").Value = 312
$ws1.Range("F39").Value = 5227
$ws1.Range("F40").Value = 1202
$ws1.Range("F41").Value = 0
$ws1.Range("F46").Value = 0
$ws1.Range("F47").Value = 1375
$ws1.Range("F48").Value = 66

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 37
$ws2.Range("F3").Value = 1
$ws2.Range("F4").Value = 16
$ws2.Range("F5").Value = 5
$ws2.Range("F6").Value = 3
$ws2.Range("F8").Value = 7
$ws2.Range("F11").Value = 1
$ws2.Range("F12").Value = 9
$ws2.Range("F13").Value = 2
$ws2.Range("F15").Value = 97
$ws2.Range("F17").Value = 7
$ws2.Range("F18").Value = 908
$ws2.Range("F19").Value = 7
$ws2.Range("F21").Value = 2

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 69
$ws4.Range("F5").Value = 27
$ws4.Range("F7").Value = 1159
$ws4.Range("F10").Value = 244
$ws4.Range("F11").Value = 7144
$ws4.Range("F13").Value = 322
$ws4.Range("F14").Value = 0
$ws4.Range("F15").Value = 3
$ws4.Range("F16").Value = 5494
$ws4.Range("F18").Value = 6267
$ws4.Range("F19").Value = 6267
$ws4.Range("F20").Value = 1096
$ws4.Range("F22").Value = 275
$ws4.Range("F23").Value = 148
$ws4.Range("F24").Value = 208
$ws4.Range("F27").Value = 10152
$ws4.Range("F28").Value = 84
$ws4.Range("F29").Value = 1929
$ws4.Range("F30").Value = 1973
$ws4.Range("F31").Value = 47
$ws4.Range("F32").Value = 2130
$ws4.Range("F34").Value = 88
$ws4.Range("F35").Value = 1033
$ws4.Range("F36").Value = 16
$ws4.Range("F37").Value = 0
$ws4.Range("F38").Value = 2087
$ws4.Range("F39").Value = 312
$ws4.Range("F40").Value = 5227
$ws4.Range("F42").Value = 667
$ws4.Range("F45").Value = 1102
$ws4.Range("F46").Value = 1078
$ws4.Range("F47").Value = 985
$ws4.Range("F48").Value = 1375
$ws4.Range("F49").Value = 66
$ws4.Range("F50").Value = 1096

